# Add a "formula" column (E) that exercises formula error handling:
#  - E2: CONCAT(...) with a valid division -> returns a string ("A1")
#  - E3: CONCAT(...) with a valid division -> returns a string ("A3")
#  - E4: CONCAT(...) with a division by zero -> #DIV/0! error
#  - E5: an (invalid) array-entered NA reference -> #NAME? error

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell E1: reuse the same formatting as the other header cells (D1)
# so it picks up the existing bordered/centered header style, then set text.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "formula"

# Data cells E2:E5: reuse the same formatting as the other data cells (A2).
$ws.Range("A2").Copy()
$ws.Range("E2:E5").PasteSpecial(-4122)

# Formulas that exercise error handling.
$ws.Range("E2").Formula = '=CONCAT("A", 3/3)'
$ws.Range("E3").Formula = '=CONCAT("A", 3/1)'
$ws.Range("E4").Formula = '=CONCAT("A", 3/0)'
$ws.Range("E5").FormulaArray = '=NA'

Write-Output "handled errors in formula cells"
